$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit removes the two data rows for account 005046919 (MARIANA,
# balance 2564) and account 005348011 (TATIANA, balance 2021.13) from the
# "Export" sheet. Locate them dynamically by account number so the edit is
# not dependent on an assumed row number, then delete both rows (shifting
# the remaining rows up), starting with the lower-most row first so the
# row index of the other target isn't invalidated.

$row1 = $ws.Cells.Find("005046919").Row
$row2 = $ws.Cells.Find("005348011").Row

if ($row1 -gt $row2) {
    $ws.Rows.Item($row1).Delete()
    $ws.Rows.Item($row2).Delete()
} else {
    $ws.Rows.Item($row2).Delete()
    $ws.Rows.Item($row1).Delete()
}
